# Apply edits per Dr Hou advice: add "sCs" group and expand LR-pairs table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Il1rn"
$ws.Cells.Item(2, 3).Value = "Il1r2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 305.026516
$ws.Cells.Item(2, 8).Value = 915.079548
$ws.Cells.Item(2, 9).Value = 0.9998851412135495
$ws.Cells.Item(2, 10).Value = 0.9998851412135495
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 265.842631
$ws.Cells.Item(2, 14).Value = 797.527893
$ws.Cells.Item(2, 15).Value = 0.99055425962745
$ws.Cells.Item(2, 16).Value = 0.99055425962745
$ws.Cells.Item(2, 17).Value = 81089.0515382036
$ws.Cells.Item(2, 18).Value = 729801.4638438324
$ws.Cells.Item(2, 19).Value = 0.9904404857672758
$ws.Cells.Item(2, 20).Value = 0.9904404857672758

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Il1rn"
$ws.Cells.Item(3, 3).Value = "Il1r2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 305.026516
$ws.Cells.Item(3, 8).Value = 915.079548
$ws.Cells.Item(3, 9).Value = 0.9998851412135495
$ws.Cells.Item(3, 10).Value = 0.9998851412135495
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.942782333333333
$ws.Cells.Item(3, 14).Value = 5.828347
$ws.Cells.Item(3, 15).Value = 0.007238986871944891
$ws.Cells.Item(3, 16).Value = 0.007238986871944891
$ws.Cells.Item(3, 17).Value = 592.6001264830173
$ws.Cells.Item(3, 18).Value = 5333.401138347156
$ws.Cells.Item(3, 19).Value = 0.007238155410697649
$ws.Cells.Item(3, 20).Value = 0.007238155410697649

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Il1rn"
$ws.Cells.Item(4, 3).Value = "Il1r2"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 305.026516
$ws.Cells.Item(4, 8).Value = 915.079548
$ws.Cells.Item(4, 9).Value = 0.9998851412135495
$ws.Cells.Item(4, 10).Value = 0.9998851412135495
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.5922433333333333
$ws.Cells.Item(4, 14).Value = 1.77673
$ws.Cells.Item(4, 15).Value = 0.002206753500604999
$ws.Cells.Item(4, 16).Value = 0.002206753500604999
$ws.Cells.Item(4, 17).Value = 180.6499205908933
$ws.Cells.Item(4, 18).Value = 1625.84928531804
$ws.Cells.Item(4, 19).Value = 0.002206500035575925
$ws.Cells.Item(4, 20).Value = 0.002206500035575925

# Row 5
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Il1rn"
$ws.Cells.Item(5, 3).Value = "Il1r2"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.035039
$ws.Cells.Item(5, 8).Value = 0.105117
$ws.Cells.Item(5, 9).Value = 0.0001148587864504919
$ws.Cells.Item(5, 10).Value = 0.0001148587864504919
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 265.842631
$ws.Cells.Item(5, 14).Value = 797.527893
$ws.Cells.Item(5, 15).Value = 0.99055425962745
$ws.Cells.Item(5, 16).Value = 0.99055425962745
$ws.Cells.Item(5, 17).Value = 9.314859947609
$ws.Cells.Item(5, 18).Value = 83.83373952848099
$ws.Cells.Item(5, 19).Value = 0.0001137738601741744
$ws.Cells.Item(5, 20).Value = 0.0001137738601741744

# Row 6
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Il1rn"
$ws.Cells.Item(6, 3).Value = "Il1r2"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.035039
$ws.Cells.Item(6, 8).Value = 0.105117
$ws.Cells.Item(6, 9).Value = 0.0001148587864504919
$ws.Cells.Item(6, 10).Value = 0.0001148587864504919
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.942782333333333
$ws.Cells.Item(6, 14).Value = 5.828347
$ws.Cells.Item(6, 15).Value = 0.007238986871944891
$ws.Cells.Item(6, 16).Value = 0.007238986871944891
$ws.Cells.Item(6, 17).Value = 0.06807315017766667
$ws.Cells.Item(6, 18).Value = 0.6126583515990001
$ws.Cells.Item(6, 19).Value = 0.0000008314612472426328
$ws.Cells.Item(6, 20).Value = 0.0000008314612472426328

# Row 7
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Il1rn"
$ws.Cells.Item(7, 3).Value = "Il1r2"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.035039
$ws.Cells.Item(7, 8).Value = 0.105117
$ws.Cells.Item(7, 9).Value = 0.0001148587864504919
$ws.Cells.Item(7, 10).Value = 0.0001148587864504919
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.5922433333333333
$ws.Cells.Item(7, 14).Value = 1.77673
$ws.Cells.Item(7, 15).Value = 0.002206753500604999
$ws.Cells.Item(7, 16).Value = 0.002206753500604999
$ws.Cells.Item(7, 17).Value = 0.02075161415666667
$ws.Cells.Item(7, 18).Value = 0.18676452741
$ws.Cells.Item(7, 19).Value = 0.0000002534650290748651
$ws.Cells.Item(7, 20).Value = 0.0000002534650290748651
